$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: Create / Input
$ws.Cells.Item(47, 1).Value = "PkRewriterApi"
$ws.Cells.Item(47, 2).Value = "PersonController"
$ws.Cells.Item(47, 3).Value = "Create"
$ws.Cells.Item(47, 5).Value = "D"
$ws.Cells.Item(47, 6).Value = "Input"
$ws.Cells.Item(47, 7).Value = '{"firstName": "Drew","lastName": "Carey"}'

# Row 48: Create / Expected
$ws.Cells.Item(48, 1).Value = "PkRewriterApi"
$ws.Cells.Item(48, 2).Value = "PersonController"
$ws.Cells.Item(48, 3).Value = "Create"
$ws.Cells.Item(48, 5).Value = "D"
$ws.Cells.Item(48, 6).Value = "Expected"
$ws.Cells.Item(48, 7).Value = '[{"id": -999001,"firstName": "Bob","lastName": "Barker"},{"id": -999002,"firstName": "Monty","lastName": "Hall"},{"id": -999301,"firstName": "Drew","lastName": "Carey"}]'

# Row 49: Create / ExpectedBypass
$ws.Cells.Item(49, 1).Value = "PkRewriterApi"
$ws.Cells.Item(49, 2).Value = "PersonController"
$ws.Cells.Item(49, 3).Value = "Create"
$ws.Cells.Item(49, 5).Value = "D"
$ws.Cells.Item(49, 6).Value = "ExpectedBypass"
$ws.Cells.Item(49, 7).Value = '[{"id": -999001,"firstName": "Bob","lastName": "Barker"},{"id": -999002,"firstName": "Monty","lastName": "Hall"},{"id": -501301,"firstName": "Drew","lastName": "Carey"}]'
